$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 778 (shifts existing rows 778-819 down to 779-820)
$ws.Rows.Item(778).Insert()

# Populate the newly inserted row with the new data point.
# The date-like text in column A must stay plain text (matching the rest of
# the column, which stores dates as literal strings, not Excel date serials).
# A leading apostrophe forces text entry, then resetting the style back to
# "Normal" clears the quote-prefix style index the text entry leaves behind.
$ws.Cells.Item(778, 1).Value = "'2026/02/04"
$ws.Cells.Item(778, 1).Style = "Normal"
$ws.Cells.Item(778, 2).Value = "水"
$ws.Cells.Item(778, 3).Value = 7
$ws.Cells.Item(778, 4).Value = 201
